# Generate Report for Handoff
# Replace the two localized files in the report with a new pair of files
# and flip their handoff/handback status/dates accordingly.

$wb = $excel.ActiveWorkbook

$oldGuid1 = "6996298e-972a-422f-874d-04f8a2062c52"
$oldGuid2 = "dd5a7bb6-83dc-49c5-b38f-53e74a0b6ef0"
$newGuid1 = "3914d24b-8b11-4bc1-8ddc-84a65dd1ee83"
$newGuid2 = "ffffb474439f-8e90-4131-8340-90fa51f53248"
$newHash  = "d97f4267a936c10b3e28a3a56e067270e7c460a3"

# ---------------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("A2").Value = "$newGuid1.md"
$ov.Range("B2").Value = "e2e\$newGuid1.md"
$ov.Range("E2").Value = "Ready for handoff"
$ov.Range("F2").Value = "Ready for handoff"
$ov.Range("G2").Value = "2016-08-18 19:07:18"

$ov.Range("A3").Value = "$newGuid2.md"
$ov.Range("B3").Value = "e2e\$newGuid2.md"
$ov.Range("E3").Value = "Ready for handoff"
$ov.Range("F3").Value = "Ready for handoff"
$ov.Range("G3").Value = "2016-08-18 19:07:18"

$ov.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid1.md"
$ov.Hyperlinks.Item(2).TextToDisplay = "e2e\$newGuid2.md"

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("A2").Value = "$newGuid1.md"
$zh.Range("C2").Value = "Ready for handoff"
$zh.Range("F2").Value = "False"
$zh.Range("G2").Value = "$newGuid1.$newHash.zh-cn.xlf"
$zh.Range("H2").Value = "2016-08-18 19:07:12"
$zh.Range("I2").Value = ""
$zh.Range("J2").Value = ""
$zh.Range("K2").Value = "0001-01-01 00:00:00"

$zh.Range("A3").Value = "$newGuid2.md"
$zh.Range("C3").Value = "Ready for handoff"
$zh.Range("F3").Value = "True"
$zh.Range("G3").Value = "$newGuid1.$newHash.zh-cn.xlf"
$zh.Range("H3").Value = "2016-08-18 19:07:12"
$zh.Range("I3").Value = ""
$zh.Range("J3").Value = ""
$zh.Range("K3").Value = "0001-01-01 00:00:00"

$zh.Hyperlinks.Item("I2").Delete()
$zh.Hyperlinks.Item("I3").Delete()
$zh.Hyperlinks.Item("A2").TextToDisplay = "$newGuid1.md"
$zh.Hyperlinks.Item("A3").TextToDisplay = "$newGuid2.md"

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")

$de.Range("A2").Value = "$newGuid1.md"
$de.Range("C2").Value = "Ready for handoff"
$de.Range("F2").Value = "False"
$de.Range("G2").Value = "$newGuid1.$newHash.de-de.xlf"
$de.Range("H2").Value = "2016-08-18 19:07:18"
$de.Range("I2").Value = ""
$de.Range("J2").Value = ""
$de.Range("K2").Value = "0001-01-01 00:00:00"

$de.Range("A3").Value = "$newGuid2.md"
$de.Range("C3").Value = "Ready for handoff"
$de.Range("F3").Value = "True"
$de.Range("G3").Value = "$newGuid1.$newHash.de-de.xlf"
$de.Range("H3").Value = "2016-08-18 19:07:18"
$de.Range("I3").Value = ""
$de.Range("J3").Value = ""
$de.Range("K3").Value = "0001-01-01 00:00:00"

$de.Hyperlinks.Item("I2").Delete()
$de.Hyperlinks.Item("I3").Delete()
$de.Hyperlinks.Item("A2").TextToDisplay = "$newGuid1.md"
$de.Hyperlinks.Item("A3").TextToDisplay = "$newGuid2.md"

# ---------------------------------------------------------------------------
# Column width adjustments
# ---------------------------------------------------------------------------
$ov.Range("E1").ColumnWidth = 17.2159881591797
$ov.Range("F1").ColumnWidth = 17.2159881591797

$zh.Range("C1").ColumnWidth = 17.2159881591797
$zh.Range("I1").ColumnWidth = 18.6506053379604
$zh.Range("J1").ColumnWidth = 21.7054770333426

$de.Range("C1").ColumnWidth = 17.2159881591797
$de.Range("I1").ColumnWidth = 18.6506053379604
$de.Range("J1").ColumnWidth = 21.7054770333426
